$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was bumped from
# 45182 (2023-09-13) to 45184 (2023-09-15) for every data row (rows 2-454).
$lastRow = 454

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}
